$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value, derived from the upstream coinranking.com scrape refresh.
$updates = @{
    "D2" = "69.397.01"
    "E2" = "  -0.44%  "
    "D3" = "3.677.45"
    "E3" = "  -0.73%  "
    "E4" = "  +0.04%  "
    "D5" = "644.97"
    "E5" = "  -5.71%  "
    "D6" = "158.49"
    "E6" = "  -1.58%  "
    "E7" = "  +0.25%  "
    "D8" = "0.498"
    "E8" = "  +0.14%  "
    "E9" = "  -1.87%  "
    "D10" = "7.08"
    "E10" = "  -1.27%  "
    "D11" = "0.445"
    "E11" = "  +0.83%  "
    "E12" = "  -1.44%  "
    "D13" = "4.301.35"
    "E13" = "  -0.71%  "
    "D14" = "32.48"
    "E14" = "  -0.31%  "
    "D15" = "3.673.84"
    "E15" = "  -0.48%  "
    "D16" = "69.431.47"
    "E16" = "  -0.29%  "
    "E17" = "  +0.38%  "
    "D18" = "15.93"
    "E18" = "  -1.00%  "
    "D19" = "6.46"
    "E19" = "  -0.31%  "
    "D20" = "466.81"
    "E20" = "  -1.08%  "
    "D21" = "9.85"
    "E21" = "  -1.12%  "
    "E22" = "  -1.39%  "
    "D23" = "79.35"
    "E23" = "  -1.33%  "
    "D24" = "3.825.85"
    "E24" = "  -0.67%  "
    "E25" = "  +0.01%  "
    "E26" = "  -1.35%  "
    "E27" = "  -2.24%  "
    "D28" = "8.96"
    "E28" = "  -3.00%  "
    "D29" = "2.63"
    "E29" = "  -3.60%  "
    "E30" = "  -3.45%  "
    "B31" = "ImmutableX"
    "C31" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
    "D31" = "1.99"
    "E31" = "  -1.59%  "
    "B32" = "Binance-PegBSC-USD"
    "C32" = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
    "D32" = "0.999"
    "E32" = "  -0.16%  "
    "D33" = "26.83"
    "E33" = "  -0.73%  "
    "E34" = "  +1.77%  "
    "E35" = "  -2.93%  "
    "D36" = "3.672.96"
    "E36" = "  -0.63%  "
    "E37" = "  -0.14%  "
    "E39" = "  -6.76%  "
    "B40" = "Monero"
    "C40" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    "D40" = "177.57"
    "E40" = "  +4.69%  "
    "B41" = "FirstDigitalUSD"
    "C41" = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
    "D41" = "1.00"
    "E41" = "  -0.13%  "
    "B42" = "Stacks"
    "C42" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "D42" = "2.22"
    "E42" = "  -4.29%  "
    "D43" = "0.0891"
    "E43" = "  -2.23%  "
    "D44" = "0.925"
    "E44" = "  -2.16%  "
    "D45" = "46.81"
    "E45" = "  -1.17%  "
    "E46" = "  -2.07%  "
    "D47" = "27.20"
    "E47" = "  -5.76%  "
    "E48" = "  -4.54%  "
    "D49" = "7.81"
    "E49" = "  -1.01%  "
    "E50" = "  -4.81%  "
    "D51" = "0.000265"
    "E51" = "  -6.54%  "
}

foreach ($addr in $updates.Keys) {
    $rng = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "7.08", "1.00")
    # keep their exact formatting instead of being coerced to a Double.
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$addr]
    # Drop back to the default style so we do not leave a stray
    # "@"-formatted style applied to the cell.
    $rng.Style = "Normal"
}
